$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original A1 ("code" header) carries the shared header style (bold, centered,
# thin border). Copy that formatting onto the new header cells (B1:F1) and onto the
# new id column (A2:A11) before we touch A1 itself, so we reuse the existing style
# entry instead of inventing a new one.
$ws.Range("A1").Copy()
$ws.Range("B1:F1").PasteSpecial(-4122)
$ws.Range("A2:A11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Old A1 ("code") no longer exists in the new layout -- the header row now starts at B1.
$ws.Range("A1").ClearContents()
$ws.Range("A1").ClearFormats()

# New header row: B1:F1
$ws.Range("B1").Value = "lang_code"
$ws.Range("C1").Value = "code"
$ws.Range("D1").Value = "name"
$ws.Range("E1").Value = "descr"
$ws.Range("F1").Value = "is_active"

$data = @(
    @(0, "eng", "POA", "Proof of Address", "Address Proof", $true),
    @(1, "fra", "POA", "Un justificatif de domicile", "Preuve dadresse", $true),
    @(2, "eng", "POI", "Proof of Identity", "Identity Proof", $true),
    @(3, "fra", "POI", "Preuve didentité", "Preuve didentité", $true),
    @(4, "eng", "POR", "Proof of Relationship", "Proof Relationship of the person", $true),
    @(5, "fra", "POR", "Preuve de relation", "Preuve de relation de la personne", $true),
    @(6, "eng", "POB", "Proof of Birth", "Proof date of birth of the person", $false),
    @(7, "fra", "POB", "Preuve de naissance", "Preuve de la date de naissance de la personne", $false),
    @(8, "eng", "POE", "Proof of Biometric Exception", "Proof of Biometric Exception", $true),
    @(9, "fra", "POE", "Preuve dexception biométrique", "Preuve dexception biométrique", $true)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}
